{"js": "// Update the multiplication problems in the practice-sheet table.\n// Each cell holds a single run like \"56\u00d752=\" which must become \"20\u00d727=\", etc.\n// We locate each old expression via Body.search (exact, case-sensitive) and\n// replace the matched range's text in place so the run's formatting\n// (TimeNewRoman, sz 30) is preserved.\n\nconst replacements = [\n  [\"56\u00d752=\", \"20\u00d727=\"],\n  [\"83\u00d784=\", \"48\u00d726=\"],\n  [\"98\u00d764=\", \"89\u00d798=\"],\n  [\"15\u00d722=\", \"99\u00d791=\"],\n  [\"28\u00d711=\", \"95\u00d791=\"],\n  [\"70\u00d774=\", \"42\u00d774=\"],\n  [\"59\u00d722=\", \"76\u00d763=\"],\n  [\"23\u00d720=\", \"82\u00d790=\"],\n  [\"63\u00d776=\", \"95\u00d782=\"],\n  [\"96\u00d719=\", \"36\u00d791=\"],\n  [\"65\u00d737=\", \"94\u00d751=\"],\n  [\"38\u00d713=\", \"78\u00d718=\"],\n  [\"71\u00d762=\", \"97\u00d738=\"],\n  [\"24\u00d719=\", \"54\u00d721=\"],\n  [\"34\u00d772=\", \"71\u00d769=\"],\n  [\"43\u00d759=\", \"53\u00d711=\"],\n  [\"76\u00d776=\", \"43\u00d798=\"],\n  [\"40\u00d783=\", \"87\u00d793=\"],\n  [\"82\u00d742=\", \"89\u00d778=\"],\n  [\"40\u00d731=\", \"74\u00d795=\"],\n  [\"85\u00d734=\", \"49\u00d754=\"],\n  [\"94\u00d799=\", \"55\u00d796=\"],\n  [\"46\u00d750=\", \"96\u00d734=\"],\n  [\"35\u00d711=\", \"84\u00d749=\"],\n  [\"94\u00d798=\", \"20\u00d790=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Could not find expression \"${oldText}\" to replace.`);\n  }\n\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update the multiplication problems in the practice-sheet table.\n# Each cell holds a single run like \"56\u00d752=\" which must become \"20\u00d727=\", etc.\n# Use Find/Replace on the whole document content so the existing run\n# formatting (TimeNewRoman, sz 30) is preserved in place.\n\n$d = $word.ActiveDocument\n\n$pairs = @(\n  @{old=\"56\u00d752=\"; new=\"20\u00d727=\"},\n  @{old=\"83\u00d784=\"; new=\"48\u00d726=\"},\n  @{old=\"98\u00d764=\"; new=\"89\u00d798=\"},\n  @{old=\"15\u00d722=\"; new=\"99\u00d791=\"},\n  @{old=\"28\u00d711=\"; new=\"95\u00d791=\"},\n  @{old=\"70\u00d774=\"; new=\"42\u00d774=\"},\n  @{old=\"59\u00d722=\"; new=\"76\u00d763=\"},\n  @{old=\"23\u00d720=\"; new=\"82\u00d790=\"},\n  @{old=\"63\u00d776=\"; new=\"95\u00d782=\"},\n  @{old=\"96\u00d719=\"; new=\"36\u00d791=\"},\n  @{old=\"65\u00d737=\"; new=\"94\u00d751=\"},\n  @{old=\"38\u00d713=\"; new=\"78\u00d718=\"},\n  @{old=\"71\u00d762=\"; new=\"97\u00d738=\"},\n  @{old=\"24\u00d719=\"; new=\"54\u00d721=\"},\n  @{old=\"34\u00d772=\"; new=\"71\u00d769=\"},\n  @{old=\"43\u00d759=\"; new=\"53\u00d711=\"},\n  @{old=\"76\u00d776=\"; new=\"43\u00d798=\"},\n  @{old=\"40\u00d783=\"; new=\"87\u00d793=\"},\n  @{old=\"82\u00d742=\"; new=\"89\u00d778=\"},\n  @{old=\"40\u00d731=\"; new=\"74\u00d795=\"},\n  @{old=\"85\u00d734=\"; new=\"49\u00d754=\"},\n  @{old=\"94\u00d799=\"; new=\"55\u00d796=\"},\n  @{old=\"46\u00d750=\"; new=\"96\u00d734=\"},\n  @{old=\"35\u00d711=\"; new=\"84\u00d749=\"},\n  @{old=\"94\u00d798=\"; new=\"20\u00d790=\"}\n)\n\nforeach ($pair in $pairs) {\n  $rng = $d.Content\n  $find = $rng.Find\n  $find.ClearFormatting()\n  $find.Replacement.ClearFormatting()\n  $result = $find.Execute($pair.old, $false, $true, $false, $false, $false, $true, 1, $false, $pair.new, 2)\n  if (-not $result) {\n    Write-Output \"WARNING: could not find '$($pair.old)'\"\n  }\n}\n"}
